$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - matches formatting/style of the other header cells (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Timestamp values for rows 2-27 (as plain text, matching column F data cells)
$timestamps = @(
    "2021-10-05 10:50:15.135886",
    "2021-10-05 10:50:15.135960",
    "2021-10-05 10:50:15.135975",
    "2021-10-05 10:50:15.135979",
    "2021-10-05 10:50:15.135983",
    "2021-10-05 10:50:15.135986",
    "2021-10-05 10:50:15.135989",
    "2021-10-05 10:50:15.135992",
    "2021-10-05 10:50:15.135996",
    "2021-10-05 10:50:15.135999",
    "2021-10-05 10:50:15.136002",
    "2021-10-05 10:50:15.136005",
    "2021-10-05 10:50:15.136008",
    "2021-10-05 10:50:15.136011",
    "2021-10-05 10:50:15.136014",
    "2021-10-05 10:50:15.136017",
    "2021-10-05 10:50:15.136021",
    "2021-10-05 10:50:15.136024",
    "2021-10-05 10:50:15.136027",
    "2021-10-05 10:50:15.136030",
    "2021-10-05 10:50:15.136035",
    "2021-10-05 10:50:15.136038",
    "2021-10-05 10:50:15.136042",
    "2021-10-05 10:50:15.136045",
    "2021-10-05 10:50:15.136048",
    "2021-10-05 10:50:15.136052"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timestamps[$i]
}
